# Generate Report for Handback
#
# Fills in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns (I/J/K) on the per-locale sheets,
# flips the "Status" column (C) from "Ready for handoff" to
# "Handed back: in sync with en-US" (which also updates the Overview
# sheet's zh-cn/de-de columns since they share the same string), and
# adds hyperlinks on the newly-populated "Latest Target File" cells.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7e4eaa4f7958505951bb014100a782f736d61a25/e2e/"

# Per-locale, per-document data: source document name, target xlf name,
# handback datetime.
$docs = @(
    @{ Name = "04d5adf0-ec1e-46c6-ae37-f077b47c4861.md" },
    @{ Name = "693b24f9-c5d0-4cdd-bedc-1177315a1cfb.md" }
)

$locales = @(
    @{
        Sheet = "zh-cn";
        Suffix = "zh-cn.xlf";
        Handoff = @{
            "04d5adf0-ec1e-46c6-ae37-f077b47c4861.md" = "2091fd8f3d5e2d97d6adce4f35ebf4dcf5a16b80";
            "693b24f9-c5d0-4cdd-bedc-1177315a1cfb.md" = "62e22bc6dc4e6c5d4757b9c2c05c90c975ef3459";
        };
        HandbackDateTime = "2016-09-03 00:32:43";
    },
    @{
        Sheet = "de-de";
        Suffix = "de-de.xlf";
        Handoff = @{
            "04d5adf0-ec1e-46c6-ae37-f077b47c4861.md" = "2091fd8f3d5e2d97d6adce4f35ebf4dcf5a16b80";
            "693b24f9-c5d0-4cdd-bedc-1177315a1cfb.md" = "62e22bc6dc4e6c5d4757b9c2c05c90c975ef3459";
        };
        HandbackDateTime = "2016-09-03 00:32:50";
    }
)

# Overview sheet: Status text lives in columns E (zh-cn) and F (de-de)
# for the two document rows.
$ws0 = $wb.Worksheets.Item("Overview")
$ws0.Range("E2").Value = $newStatus
$ws0.Range("F2").Value = $newStatus
$ws0.Range("E3").Value = $newStatus
$ws0.Range("F3").Value = $newStatus

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    $row = 2
    foreach ($doc in $docs) {
        $docName = $doc.Name
        $docStem = $docName.Substring(0, $docName.Length - 3)  # strip ".md"
        $hash = $locale.Handoff[$docName]
        $targetFile = "$docStem.$hash.$($locale.Suffix)"

        # Status
        $ws.Range("C$row").Value = $newStatus

        # Latest Target File / Latest Handback File (same generated xlf)
        $ws.Range("I$row").Value = $docName
        $ws.Range("J$row").Value = $targetFile

        # Latest Handback DateTime
        $ws.Range("K$row").Value = $locale.HandbackDateTime

        # Add a hyperlink on the newly-populated target-file cell (I),
        # mirroring the one already on the source file name (A).
        $ws.Hyperlinks.Add($ws.Range("I$row"), $repoBase + $docName, "", "", $docName) | Out-Null

        $row = $row + 1
    }

    # Widen the newly-populated columns to fit their longer content.
    $ws.Columns.Item(3).ColumnWidth = 29.98
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
}

# Overview sheet: zh-cn / de-de columns also grew to fit the new status text.
$ws0.Columns.Item(5).ColumnWidth = 29.98
$ws0.Columns.Item(6).ColumnWidth = 29.98
